$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.57920241355896
$ws.Range("B1").Value = 2.263180494308472
$ws.Range("C1").Value = 6.248862266540527
$ws.Range("D1").Value = 1.579378008842468
$ws.Range("E1").Value = 0.9066450595855713
